$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUPPLY_ASSEMBLIES")

# Header G1
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "type_el_pv"

# Data G2:G18
for ($row = 2; $row -le 18; $row++) {
    $c = $ws.Cells.Item($row, 7)
    $c.Value = "SUPPLY_ELECTRICITY_PV_AS0"
    $c.Interior.Pattern = 1
    $c.Interior.PatternColorIndex = -4105
    $c.Interior.ThemeColor = 0
    $c.Interior.TintAndShade = 0
    $c.Borders.LineStyle = 1
}

$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

$ws.Activate()
$ws.Range("G1").Select() | Out-Null
